# Edit script: "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Semantic changes applied to Hoja1:
#   1. VALOR MORA total (E11) updated:        71172   -> 142344
#   2. Cant. Periodos (F13) updated:           1       -> 2
#   3. A new worker/period detail row is inserted right below the existing
#      one (old row 16 -> stays row 16; new row 17 added), duplicating the
#      same worker line but for period 2509 instead of 2508.
#   4. Everything below (the two signature blocks) shifts down by one row
#      as a natural consequence of the row insertion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new detail row under row 16, copying its formatting ------
# Using "copy whole row, then insert copied cells" reproduces Excel's
# native "Insert Copied Cells" behaviour: values/format come along for the
# ride and every row/mergeCell reference below the insertion point shifts
# down by one automatically.
$ws.Rows("16:16").Copy()
$ws.Rows("17:17").Insert(-4121)  # -4121 = xlShiftDown
$excel.CutCopyMode = 0

# Re-apply row 16's exact cell formatting (borders included) onto the new
# row 17 so the duplicated line keeps the same boxed look as the original.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 2. Update the new row's period (E17) to the new period "2509" -------
$ws.Range("E17").Value = "2509"

# --- 3. Update summary fields --------------------------------------------
$ws.Range("E11").Value = 142344
$ws.Range("F13").Value = 2
